$wb = $excel.ActiveWorkbook

# Rename sheets: GeneralTaxRateWeekly -> GeneralTaxRateMonthly ; ProcessPayrollForWeeklyTax -> ProcessPayrollForMonthlyTax
$wsGeneral = $wb.Worksheets.Item("GeneralTaxRateWeekly")
$wsGeneral.Name = "GeneralTaxRateMonthly"

$wsProcess = $wb.Worksheets.Item("ProcessPayrollForWeeklyTax")
$wsProcess.Name = "ProcessPayrollForMonthlyTax"

# Update "first" sheet references to the renamed sheets
$wsFirst = $wb.Worksheets.Item("first")
$wsFirst.Range("A3").Value = "GeneralTaxRateMonthly"
$wsFirst.Range("A4").Value = "ProcessPayrollForMonthlyTax"

# Update "DO NOT TOUCH AUTOMATION EMP 107" -> "DO NOT TOUCH AUTOMATION EMP 105" on the relevant sheets
$wsGeneral.Range("A2").Value = "DO NOT TOUCH AUTOMATION EMP 105"
$wsProcess.Range("B2").Value = "DO NOT TOUCH AUTOMATION EMP 105"

$wsReports = $wb.Worksheets.Item("TestReports")
$wsReports.Range("B2").Value = "DO NOT TOUCH AUTOMATION EMP 105"

# Update selections / active sheet
$wsFirst.Range("F5").Select()
$wsGeneral.Range("G10").Select()
$wsProcess.Range("H13").Select()
$wsReports.Range("M4").Select()

$wsGeneral.Activate()
